$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.743
$ws.Range("B4").Value = 6.612
$ws.Range("B7").Value = 6.866
$ws.Range("B8").Value = 6.444
$ws.Range("A11").Value = -21.584
$ws.Range("A12").Value = -21.36
$ws.Range("B12").Value = 6.695
$ws.Range("B14").Value = 6.567
$ws.Range("A15").Value = -21.098
$ws.Range("B22").Value = 6.964
